# Generate Report for Handback
# - Updates "Ready for handoff" status to "Handed back: in sync with en-US"
#   on the Overview sheet and each per-language sheet.
# - Fills in the "Latest Target File" / "Latest Handback File" columns (E/F)
#   and the "Latest Handback DateTime" (G) for the two real source files on
#   each per-language sheet, now that handback has happened.
# - Rebuilds the hyperlinks for each per-language sheet so the newly
#   populated E/F cells link to the same targets as the corresponding
#   A/C cells.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: just the status text changes (B2:C2, B3:C3)
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = $newStatus
$wsOverview.Range("C2").Value = $newStatus
$wsOverview.Range("B3").Value = $newStatus
$wsOverview.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# Per-language sheets
# ---------------------------------------------------------------------
$langInfo = @{
    "zh-cn" = @{
        HandbackDateTime = "2016-03-10 12:10:01"
        XlfName          = "42c2c268-8452-4b8e-9a61-a8cbab82c5ed.49ad32d97d9a23db737c8e9dcd7920e9bb3249f9.zh-cn.xlf"
        XlfUrl           = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/0b080f3de8b8eaa47c55de246bfff8a1f1666db0/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/42c2c268-8452-4b8e-9a61-a8cbab82c5ed.49ad32d97d9a23db737c8e9dcd7920e9bb3249f9.zh-cn.xlf"
    }
    "de-de" = @{
        HandbackDateTime = "2016-03-10 12:10:10"
        XlfName          = "42c2c268-8452-4b8e-9a61-a8cbab82c5ed.49ad32d97d9a23db737c8e9dcd7920e9bb3249f9.de-de.xlf"
        XlfUrl           = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/2577ba9795aae7d2295e9d78917629f2470e4689/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/42c2c268-8452-4b8e-9a61-a8cbab82c5ed.49ad32d97d9a23db737c8e9dcd7920e9bb3249f9.de-de.xlf"
    }
}

$mdName = "42c2c268-8452-4b8e-9a61-a8cbab82c5ed.md"
$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/14e84e181fcdf7d6c786f873d55a4f03dae2ab8b/e2e/42c2c268-8452-4b8e-9a61-a8cbab82c5ed.md"
$ffffName = "ffff181c2732-a646-4607-9337-a000a24fdebc.md"
$ffffUrl = "https://github.com/OpenLocalizationTest/oltest/blob/14e84e181fcdf7d6c786f873d55a4f03dae2ab8b/e2e/ffff181c2732-a646-4607-9337-a000a24fdebc.md"
$configName = ".localization-config"
$configUrl = "https://github.com/OpenLocalizationTest/oltest/blob/14e84e181fcdf7d6c786f873d55a4f03dae2ab8b/.localization-config"

$langNames = @("zh-cn", "de-de")

foreach ($lang in $langNames) {
    $info = $langInfo[$lang]
    $ws = $wb.Worksheets.Item($lang)

    # --- status text (row 2 & 3) ---
    $ws.Range("B2").Value = $newStatus
    $ws.Range("B3").Value = $newStatus

    # --- newly-populated "Latest Target File" / "Latest Handback File" ---
    $ws.Range("E2").Value = $mdName
    $ws.Range("F2").Value = $info.XlfName
    $ws.Range("E3").Value = $mdName
    $ws.Range("F3").Value = $info.XlfName

    # --- "Latest Handback DateTime" now populated for both data rows ---
    $ws.Range("G2").Value = $info.HandbackDateTime
    $ws.Range("G3").Value = $info.HandbackDateTime

    # --- rebuild hyperlinks in the canonical row-major / column-major order
    #     (A2, C2, E2, F2, A3, C3, E3, F3, A4) ---
    $ws.Range("A1").Hyperlinks.Delete()

    $ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl, $null, $null, $mdName) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C2"), $info.XlfUrl, $null, $null, $info.XlfName) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("E2"), $mdUrl, $null, $null, $mdName) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F2"), $info.XlfUrl, $null, $null, $info.XlfName) | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A3"), $ffffUrl, $null, $null, $ffffName) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("C3"), $info.XlfUrl, $null, $null, $info.XlfName) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("E3"), $mdUrl, $null, $null, $mdName) | Out-Null
    $ws.Hyperlinks.Add($ws.Range("F3"), $info.XlfUrl, $null, $null, $info.XlfName) | Out-Null

    $ws.Hyperlinks.Add($ws.Range("A4"), $configUrl, $null, $null, $configName) | Out-Null
}
